$wb = $excel.ActiveWorkbook
$aw = $excel.ActiveWindow
$aw.Left = -120
$aw.Top = -120
$aw.Width = 29040
$aw.Height = 13920
Write-Host $aw.Left
Write-Host $aw.Top
Write-Host $aw.Width
Write-Host $aw.Height
